# Update TPM-derived NATMI edge statistics for Lama2-Itga7 (OldD7) per new TPM values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.667069666666666
$ws.Range("H2").Value = 11.001209
$ws.Range("I2").Value = 0.01298011522000835
$ws.Range("J2").Value = 0.01298011522000835
$ws.Range("M2").Value = 3.691731666666667
$ws.Range("N2").Value = 11.075195
$ws.Range("O2").Value = 0.04949987503010053
$ws.Range("P2").Value = 0.04949987503010053
$ws.Range("Q2").Value = 13.53783721230611
$ws.Range("R2").Value = 121.840534910755
$ws.Range("S2").Value = 0.0006425140812667193
$ws.Range("T2").Value = 0.0006425140812667194

# Row 3
$ws.Range("G3").Value = 3.667069666666666
$ws.Range("H3").Value = 11.001209
$ws.Range("I3").Value = 0.01298011522000835
$ws.Range("J3").Value = 0.01298011522000835
$ws.Range("O3").Value = 0.05430547069958891
$ws.Range("P3").Value = 0.05430547069958892
$ws.Range("Q3").Value = 14.85213087147466
$ws.Range("R3").Value = 133.669177843272
$ws.Range("S3").Value = 0.0007048912667574517
$ws.Range("T3").Value = 0.0007048912667574518

# Row 4
$ws.Range("G4").Value = 3.667069666666666
$ws.Range("H4").Value = 11.001209
$ws.Range("I4").Value = 0.01298011522000835
$ws.Range("J4").Value = 0.01298011522000835
$ws.Range("M4").Value = 66.272152
$ws.Range("N4").Value = 198.816456
$ws.Range("O4").Value = 0.8885974220704449
$ws.Range("P4").Value = 0.888597422070445
$ws.Range("Q4").Value = 243.0245983439227
$ws.Range("R4").Value = 2187.221385095304
$ws.Range("S4").Value = 0.01153409692267677
$ws.Range("T4").Value = 0.01153409692267677

# Row 5
$ws.Range("G5").Value = 3.667069666666666
$ws.Range("H5").Value = 11.001209
$ws.Range("I5").Value = 0.01298011522000835
$ws.Range("J5").Value = 0.01298011522000835
$ws.Range("M5").Value = 0.5666063333333334
$ws.Range("N5").Value = 1.699819
$ws.Range("O5").Value = 0.007597232199865597
$ws.Range("P5").Value = 0.007597232199865597
$ws.Range("Q5").Value = 2.077784897907889
$ws.Range("R5").Value = 18.700064081171
$ws.Range("S5").Value = 0.00009861294930741298
$ws.Range("T5").Value = 0.00009861294930741299

# Row 6
$ws.Range("I6").Value = 0.5954329572989919
$ws.Range("J6").Value = 0.595432957298992
$ws.Range("M6").Value = 3.691731666666667
$ws.Range("N6").Value = 11.075195
$ws.Range("O6").Value = 0.04949987503010053
$ws.Range("P6").Value = 0.04949987503010053
$ws.Range("Q6").Value = 621.0171720456101
$ws.Range("R6").Value = 5589.154548410491
$ws.Range("S6").Value = 0.02947385697510329
$ws.Range("T6").Value = 0.02947385697510329

# Row 7
$ws.Range("I7").Value = 0.5954329572989919
$ws.Range("J7").Value = 0.595432957298992
$ws.Range("O7").Value = 0.05430547069958891
$ws.Range("P7").Value = 0.05430547069958892
$ws.Range("S7").Value = 0.03233526701616998
$ws.Range("T7").Value = 0.03233526701616999

# Row 8
$ws.Range("I8").Value = 0.5954329572989919
$ws.Range("J8").Value = 0.595432957298992
$ws.Range("M8").Value = 66.272152
$ws.Range("N8").Value = 198.816456
$ws.Range("O8").Value = 0.8885974220704449
$ws.Range("P8").Value = 0.888597422070445
$ws.Range("Q8").Value = 11148.19497636389
$ws.Range("R8").Value = 100333.754787275
$ws.Range("S8").Value = 0.5291001908716655
$ws.Range("T8").Value = 0.5291001908716657

# Row 9
$ws.Range("I9").Value = 0.5954329572989919
$ws.Range("J9").Value = 0.595432957298992
$ws.Range("M9").Value = 0.5666063333333334
$ws.Range("N9").Value = 1.699819
$ws.Range("O9").Value = 0.007597232199865597
$ws.Range("P9").Value = 0.007597232199865597
$ws.Range("Q9").Value = 95.31360742356203
$ws.Range("R9").Value = 857.8224668120581
$ws.Range("S9").Value = 0.004523642436053099
$ws.Range("T9").Value = 0.004523642436053099

# Row 10
$ws.Range("G10").Value = 110.4727123333333
$ws.Range("H10").Value = 331.418137
$ws.Range("I10").Value = 0.3910338949346852
$ws.Range("J10").Value = 0.3910338949346853
$ws.Range("M10").Value = 3.691731666666667
$ws.Range("N10").Value = 11.075195
$ws.Range("O10").Value = 0.04949987503010053
$ws.Range("P10").Value = 0.04949987503010053
$ws.Range("Q10").Value = 407.8356104235239
$ws.Range("R10").Value = 3670.520493811715
$ws.Range("S10").Value = 0.01935612893180038
$ws.Range("T10").Value = 0.01935612893180038

# Row 11
$ws.Range("G11").Value = 110.4727123333333
$ws.Range("H11").Value = 331.418137
$ws.Range("I11").Value = 0.3910338949346852
$ws.Range("J11").Value = 0.3910338949346853
$ws.Range("O11").Value = 0.05430547069958891
$ws.Range("P11").Value = 0.05430547069958892
$ws.Range("Q11").Value = 447.4295092388772
$ws.Range("R11").Value = 4026.865583149895
$ws.Range("S11").Value = 0.02123527972392168
$ws.Range("T11").Value = 0.02123527972392168

# Row 12
$ws.Range("G12").Value = 110.4727123333333
$ws.Range("H12").Value = 331.418137
$ws.Range("I12").Value = 0.3910338949346852
$ws.Range("J12").Value = 0.3910338949346853
$ws.Range("M12").Value = 66.272152
$ws.Range("N12").Value = 198.816456
$ws.Range("O12").Value = 0.8885974220704449
$ws.Range("P12").Value = 0.888597422070445
$ws.Range("Q12").Value = 7321.264383606942
$ws.Range("R12").Value = 65891.37945246248
$ws.Range("S12").Value = 0.3474717109811265
$ws.Range("T12").Value = 0.3474717109811266

# Row 13
$ws.Range("G13").Value = 110.4727123333333
$ws.Range("H13").Value = 331.418137
$ws.Range("I13").Value = 0.3910338949346852
$ws.Range("J13").Value = 0.3910338949346853
$ws.Range("M13").Value = 0.5666063333333334
$ws.Range("N13").Value = 1.699819
$ws.Range("O13").Value = 0.007597232199865597
$ws.Range("P13").Value = 0.007597232199865597
$ws.Range("Q13").Value = 62.59453846857812
$ws.Range("R13").Value = 563.3508462172031
$ws.Range("S13").Value = 0.002970775297836651
$ws.Range("T13").Value = 0.002970775297836652

# Row 14
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.1562396666666667
$ws.Range("H14").Value = 0.468719
$ws.Range("I14").Value = 0.0005530325463144183
$ws.Range("J14").Value = 0.0005530325463144184
$ws.Range("M14").Value = 3.691731666666667
$ws.Range("N14").Value = 11.075195
$ws.Range("O14").Value = 0.04949987503010053
$ws.Range("P14").Value = 0.04949987503010053
$ws.Range("Q14").Value = 0.5767949250227778
$ws.Range("R14").Value = 5.191154325205001
$ws.Range("S14").Value = 0.00002737504193014199
$ws.Range("T14").Value = 0.00002737504193014199

# Row 15
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.1562396666666667
$ws.Range("H15").Value = 0.468719
$ws.Range("I15").Value = 0.0005530325463144183
$ws.Range("J15").Value = 0.0005530325463144184
$ws.Range("O15").Value = 0.05430547069958891
$ws.Range("P15").Value = 0.05430547069958892
$ws.Range("Q15").Value = 0.6327918985946666
$ws.Range("R15").Value = 5.695127087352
$ws.Range("S15").Value = 0.00003003269273979669
$ws.Range("T15").Value = 0.0000300326927397967

# Row 16
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.1562396666666667
$ws.Range("H16").Value = 0.468719
$ws.Range("I16").Value = 0.0005530325463144183
$ws.Range("J16").Value = 0.0005530325463144184
$ws.Range("M16").Value = 66.272152
$ws.Range("N16").Value = 198.816456
$ws.Range("O16").Value = 0.8885974220704449
$ws.Range("P16").Value = 0.888597422070445
$ws.Range("Q16").Value = 10.35433893776267
$ws.Range("R16").Value = 93.189050439864
$ws.Range("S16").Value = 0.000491423294976046
$ws.Range("T16").Value = 0.0004914232949760461

# Row 17
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 0.1562396666666667
$ws.Range("H17").Value = 0.468719
$ws.Range("I17").Value = 0.0005530325463144183
$ws.Range("J17").Value = 0.0005530325463144184
$ws.Range("M17").Value = 0.5666063333333334
$ws.Range("N17").Value = 1.699819
$ws.Range("O17").Value = 0.007597232199865597
$ws.Range("P17").Value = 0.007597232199865597
$ws.Range("Q17").Value = 0.08852638465122224
$ws.Range("R17").Value = 0.796737461861
$ws.Range("S17").Value = 0.00000420151666843356
$ws.Range("T17").Value = 0.000004201516668433561
